# Generate Report for Handback
# Update the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamps for the 8bc13b2a... row on each sheet.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: row 4 (8bc13b2a...) - Latest HO Xliff Generate Date
$overview.Range("G4").Value = "2016-11-09 00:21:30"

# zh-cn sheet: row 4 (8bc13b2a...) - Correspond Handoff Datetime / Correspond Handback DateTime
$zhcn.Range("H4").Value = "2016-11-09 00:21:16"
$zhcn.Range("K4").Value = "2016-11-09 00:22:12"

# de-de sheet: row 4 (8bc13b2a...) - Correspond Handoff Datetime / Correspond Handback DateTime
$dede.Range("H4").Value = "2016-11-09 00:21:30"
$dede.Range("K4").Value = "2016-11-09 00:22:30"
